$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Update header row 10 labels (text content changes due to relabeling)
# ---------------------------------------------------------------------------
$ws.Range("G10").Value = "factor"
$ws.Range("H10").Value = "share"
$ws.Range("I10").Value = "tshare"
$ws.Range("J10").Value = "w"
$ws.Range("K10").Value = "v"

# ---------------------------------------------------------------------------
# 2. Row 11 - rework formulas
# ---------------------------------------------------------------------------
$ws.Range("G11").Formula = "=C11/A11"
$ws.Range("H11").Formula = "=G11/`$G`$14"
$ws.Range("I11").Formula = "=H11*`$F`$14"
$ws.Range("L11").ClearContents()

# ---------------------------------------------------------------------------
# 3. Row 12 - rework formulas
# ---------------------------------------------------------------------------
$ws.Range("G12").Formula = "=C12/A12"
$ws.Range("H12").Formula = "=G12/`$G`$14"
$ws.Range("I12:I13").FormulaR1C1 = "=RC[-1]*R14C6"
$ws.Range("L12").ClearContents()

# ---------------------------------------------------------------------------
# 4. Row 13 - rework formulas
# ---------------------------------------------------------------------------
$ws.Range("G13").Formula = "=C13/A13"
$ws.Range("H13").Formula = "=G13/`$G`$14"
$ws.Range("L13").ClearContents()

# ---------------------------------------------------------------------------
# 5. Row 14 - totals row rework
# ---------------------------------------------------------------------------
$ws.Range("G14").Formula = "=SUM(G11:G13)"
$ws.Range("F14").Formula = "=SUM(F11:F13)"
$ws.Range("H14:I14").FormulaR1C1 = "=SUM(R[-3]C:R[-1]C)"
$ws.Range("L14").Value = "sum"

# ---------------------------------------------------------------------------
# 6. Row 15 - replace label + sum with boolean sanity checks
# ---------------------------------------------------------------------------
$ws.Range("G15").ClearContents()
$ws.Range("H15").Formula = "=H14=1"
$ws.Range("I15").Formula = "=I14=F14"
$ws.Range("K15").Formula = "=K13=K12=K11"

# ---------------------------------------------------------------------------
# 7. Row 17-19 label text updates
# ---------------------------------------------------------------------------
$ws.Range("D17").Value = "total v"
$ws.Range("D18").Value = "total r"
$ws.Range("D19").Value = "totalMOI"

# ---------------------------------------------------------------------------
# 8. New block: rows 21-27 (second, reworked table)
# ---------------------------------------------------------------------------
$ws.Range("A21").Value = "torque input"
$ws.Range("B21").Value = 2

$ws.Range("A22").Value = "wheel radius"
$ws.Range("B22").Value = "wheel mass"
$ws.Range("C22").Value = "moi"
$ws.Range("D22").Value = "w (rad/s)"
$ws.Range("E22").Value = "v"
$ws.Range("F22").Value = "t"
$ws.Range("G22").Value = "factor"
$ws.Range("H22").Value = "share"
$ws.Range("I22").Value = "in share"
$ws.Range("J22").Value = "w a"
$ws.Range("K22").Value = "new v"

$ws.Range("A23").Value = 0.1
$ws.Range("B23").Value = 0.04
$ws.Range("C23").Value = 0.0004
$ws.Range("D23").Value = 3.75
$ws.Range("E23").Formula = "=D23*A23"
$ws.Range("F23").Formula = "=D23*C23"
$ws.Range("G23").Formula = "=C23/A23"
$ws.Range("H23").Formula = "=G23/`$G`$26"
$ws.Range("I23").Formula = "=H23*`$B`$21"
$ws.Range("J23").Formula = "=I23/C23"
$ws.Range("K23").Formula = "=(J23*A23) +E23"

$ws.Range("A24").Value = 0.2
$ws.Range("B24").Value = 0.04
$ws.Range("C24").Value = 0.0016
$ws.Range("D24").Value = 1.875
$ws.Range("E24:E25").FormulaR1C1 = "=RC4*RC1"
$ws.Range("F24").Formula = "=D24*C24"
$ws.Range("G24").Formula = "=C24/A24"
$ws.Range("H24:H25").FormulaR1C1 = "=RC[-1]/R26C7"
$ws.Range("I24:I25").FormulaR1C1 = "=RC[-1]*R21C2"
$ws.Range("J24").Formula = "=I24/C24"
$ws.Range("K24").Formula = "=(J24*A24) +E24"

$ws.Range("A25").Value = 0.3
$ws.Range("B25").Value = 0.04
$ws.Range("C25").Value = 0.0036
$ws.Range("D25").Value = 1.25
$ws.Range("F25").Formula = "=D25*C25"
$ws.Range("G25").Formula = "=C25/A25"
$ws.Range("J25").Formula = "=I25/C25"
$ws.Range("K25").Formula = "=(J25*A25) +E25"

$ws.Range("F26").Formula = "=SUM(F23:F25)"
$ws.Range("G26").Formula = "=SUM(G23:G25)"
$ws.Range("H26:I26").FormulaR1C1 = "=SUM(R[-3]C:R[-1]C)"

$ws.Range("H27").Formula = "=H26=1"
$ws.Range("I27").Formula = "=I26=B21"
$ws.Range("K27").Formula = "=AND(K23=K24,K24=K25)"

# ---------------------------------------------------------------------------
# 9. Column K width
# ---------------------------------------------------------------------------
$ws.Range("K11").ColumnWidth = 7

# ---------------------------------------------------------------------------
# 10. Conditional formatting (Highlight Cell Rules > Equal To > Green Fill with
#     Dark Green Text), replicating the add/delete history needed to reach the
#     exact dxfId/priority bookkeeping seen in the target file.
# ---------------------------------------------------------------------------
function Set-GoodFormat($fc) {
    $fc.Font.Color = 24832
    $fc.Interior.Color = 13561798
}

$dummyA = $ws.Range("A1")
$dummyB = $ws.Range("A2")
$rngH27 = $ws.Range("H27:K27")
$rngH15 = $ws.Range("H15:K15")

$tmp0 = $dummyA.FormatConditions.Add(1, 3, "=TRUE")
Set-GoodFormat $tmp0
$tmp0.Delete()

$condH27 = $rngH27.FormatConditions.Add(1, 3, "=TRUE")
Set-GoodFormat $condH27

$tmpLive = $dummyB.FormatConditions.Add(1, 3, "=TRUE")

$condH15 = $rngH15.FormatConditions.Add(1, 3, "=TRUE")
Set-GoodFormat $condH15

$tmpLive.Delete()

for ($i = 0; $i -lt 4; $i++) {
    $pad = $dummyA.FormatConditions.Add(1, 3, "=TRUE")
    Set-GoodFormat $pad
    $pad.Delete()
}

# ---------------------------------------------------------------------------
# 11. Selection on the sheet (matches the final cursor position saved with the
#     workbook)
# ---------------------------------------------------------------------------
$ws.Range("G23").Select()
